$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: the empty paragraph right before "Build has to be:" gets two tab
# runs + "reactphysics3d.lib (server)" text appended.
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "`r") {
        $nxt = $p.Next()
        if ($nxt -ne $null -and $nxt.Range.Text.StartsWith("Build has to be:")) {
            $targetPara = $p
        }
    }
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' +
      '<w:p w:rsidR="00CA537E" w:rsidRDefault="00CA537E" w:rsidP="00AB2365">' +
      '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>reactphysics3d.lib (server)</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    $r2 = $r.Duplicate
    $r2.Collapse(0)
    $r2.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Part 2: "Build has to be: x86" paragraph.
#   - "x86" run gets Italic + single Underline.
#   - " " (plain) + bold "!!!!!! Else you will get lots of errors that it
#     can" are inserted right after "x86" (before the _GoBack bookmark).
#   - bold "not find the libraries." is inserted at the very end of the
#     paragraph (after the bookmark).
# ---------------------------------------------------------------------------

# Insert the plain space before the bookmark.
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore(" ")

# Insert the bold warning text before the bookmark (re-fetch, position moved).
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore("!!!!!! Else you will get lots of errors that it can")

# Format that just-inserted warning text as bold via Find (unique text).
$fr = $d.Content
$fr.Find.Execute("!!!!!! Else you will get lots of errors that it can", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fr.Bold = 1

# Now style "x86" (must happen AFTER the inserts above so the new runs don't
# inherit the italic/underline formatting).
$xr = $d.Content
$xr.Find.Execute("x86", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xr.Italic = 1
$xr.Font.Underline = 1

# Append the final bold run at the very end of the paragraph (after the
# bookmark).
$p19 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Build has to be:")) {
        $p19 = $p
    }
}
$pr = $p19.Range
$pr.Collapse(0)
$pr.InsertAfter("not find the libraries.")

$lr = $d.Content
$lr.Find.Execute("not find the libraries.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lr.Bold = 1
$lr.LanguageID = "en-US"

Write-Host "FINAL: [$($p19.Range.Text)]"
